# Generate Report for Handback
# c5306d91-0bb6-470b-8583-15d1d98f19a9.md has been handed back (in both the
# zh-cn and de-de target languages). Update the per-language status tables
# and the Overview roll-up to reflect the handback, and record the new
# "Latest Handback" file/datetime plus the hyperlink on the newly-populated
# "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$sourceFile       = "c5306d91-0bb6-470b-8583-15d1d98f19a9.md"
$sourceUrl        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a09b0ff13d0c2d5ddc12ffc6c51804678d68b18/e2e/c5306d91-0bb6-470b-8583-15d1d98f19a9.md"

# --- zh-cn -----------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhHandbackFile = "c5306d91-0bb6-470b-8583-15d1d98f19a9.c7b8ca099646ae333dacf5eefdc16ca675f61dae.zh-cn.xlf"
$zhHandbackTime = "2016-09-07 05:55:34"

foreach ($r in 3, 4) {
    $zh.Range("C$r").Value = $statusHandedBack
    $zh.Range("I$r").Value = $sourceFile
    $zh.Range("J$r").Value = $zhHandbackFile
    $zh.Range("K$r").Value = $zhHandbackTime
}
$zh.Hyperlinks.Add($zh.Range("I3"), $sourceUrl, "", "", $sourceFile)
$zh.Hyperlinks.Add($zh.Range("I4"), $sourceUrl, "", "", $sourceFile)

# --- de-de -------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deHandbackFile = "c5306d91-0bb6-470b-8583-15d1d98f19a9.c7b8ca099646ae333dacf5eefdc16ca675f61dae.de-de.xlf"
$deHandbackTime = "2016-09-07 05:55:51"

foreach ($r in 3, 4) {
    $de.Range("C$r").Value = $statusHandedBack
    $de.Range("I$r").Value = $sourceFile
    $de.Range("J$r").Value = $deHandbackFile
    $de.Range("K$r").Value = $deHandbackTime
}
$de.Hyperlinks.Add($de.Range("I3"), $sourceUrl, "", "", $sourceFile)
$de.Hyperlinks.Add($de.Range("I4"), $sourceUrl, "", "", $sourceFile)

# --- Overview roll-up ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
foreach ($r in 3, 4) {
    $ov.Range("E$r").Value = $statusHandedBack
    $ov.Range("F$r").Value = $statusHandedBack
}
